# Apply the edit described by the commit "Responded to the first post."
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week1")

# The author posted a response to Discussion question 1 (row 7),
# recording 10 minutes (as a day fraction) in the "Actual time length
# to complete" column (C).
$ws.Range("C7").Value = 10 / 1440

# Move the active selection to the cell that was just edited.
$ws.Range("C8").Select()

# Reflect the updated on-screen window position recorded for the workbook.
$excel.ActiveWindow.Left = -27820
$excel.ActiveWindow.Top = 760
